$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.888.81"
$ws.Range("E2").Value = "  +6.61%  "
# Row 3
$ws.Range("D3").Value = "3.628.72"
$ws.Range("E3").Value = "  +4.28%  "
# Row 4
$ws.Range("E4").Value = "  -0.07%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "420.15"
$ws.Range("E5").Value = "  +1.23%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.49"
$ws.Range("E6").Value = "  +0.24%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.653"
$ws.Range("E7").Value = "  +2.34%  "
# Row 8
$ws.Range("D8").Value = "3.619.98"
$ws.Range("E8").Value = "  +4.27%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.07%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.764"
$ws.Range("E10").Value = "  +1.71%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.196"
$ws.Range("E11").Value = "  +25.25%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000434"
$ws.Range("E12").Value = "  +88.63%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.96"
$ws.Range("E13").Value = "  -1.15%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.82"
$ws.Range("E14").Value = "  +1.49%  "
# Row 15
$ws.Range("D15").Value = "4.200.79"
$ws.Range("E15").Value = "  +4.23%  "
# Row 16
$ws.Range("E16").Value = "  +0.24%  "
# Row 17
$ws.Range("D17").Value = "3.681.50"
$ws.Range("E17").Value = "  +6.08%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.04"
$ws.Range("E18").Value = "  -0.99%  "
# Row 19
$ws.Range("E19").Value = "  +2.26%  "
# Row 20
$ws.Range("D20").Value = "67.704.03"
$ws.Range("E20").Value = "  +6.65%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.40"
$ws.Range("E21").Value = "  +0.16%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "459.23"
$ws.Range("E22").Value = "  +0.04%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.02"
$ws.Range("E23").Value = "  -1.15%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.37"
$ws.Range("E24").Value = "  +1.34%  "
# Row 25
$ws.Range("E25").Value = "  -6.72%  "
# Row 26
$ws.Range("E26").Value = "  -1.18%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "35.53"
$ws.Range("E27").Value = "  +6.52%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.25"
$ws.Range("E28").Value = "  -1.77%  "
# Row 29
$ws.Range("E29").Value = "  +3.76%  "
# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.79"
$ws.Range("E30").Value = "  +4.87%  "
# Row 31
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.21"
$ws.Range("E31").Value = "  -2.65%  "
# Row 32
$ws.Range("E32").Value = "  +5.45%  "
# Row 33
$ws.Range("E33").Value = "  -4.12%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.158"
$ws.Range("E34").Value = "  -7.27%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "39.94"
$ws.Range("E35").Value = "  +0.25%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.15%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.08"
$ws.Range("E37").Value = "  -2.62%  "
# Row 38
$ws.Range("D38").Value = "0.0₃0797"
$ws.Range("E38").Value = "  +23.17%  "
# Row 39
$ws.Range("E39").Value = "  +0.84%  "
# Row 40
$ws.Range("E40").Value = "  +8.61%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.13%  "
# Row 42
$ws.Range("E42").Value = "  -2.89%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "147.88"
$ws.Range("E43").Value = "  +0.78%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.92"
$ws.Range("E44").Value = "  -4.87%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.24"
$ws.Range("E45").Value = "  -2.28%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.62"
$ws.Range("E46").Value = "  +12.08%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.26"
$ws.Range("E47").Value = "  -5.59%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.168"
$ws.Range("E48").Value = "  +19.81%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.301"
$ws.Range("E49").Value = "  -4.51%  "
# Row 51
$ws.Range("E51").Value = "  +12.69%  "
